$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new year column (P) for 2022, matching the
# formatting of the preceding column (O) for rows 3-5.
$ws.Range("O3:O5").Copy() | Out-Null
$ws.Range("P3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# New header value for the added year column
$ws.Range("P4").Value = 2022

# Updated data series values (row 5) including the new 2022 data point
$ws.Range("M5").Value = 2.6
$ws.Range("N5").Value = 2.4
$ws.Range("O5").Value = 3.3
$ws.Range("P5").Value = 2.6

# Move the active selection to the newly added cell
$ws.Range("P3").Select() | Out-Null
